# hotfix double-encoding-issue by using triple-braces
#
# The template placeholders in the shared-string table use Handlebars-style
# "{{ }}" delimiters. Because the downstream templating step HTML/URL-encodes
# the merge output, double braces get double-encoded. Switching every
# placeholder to triple braces ("{{{ }}}") fixes that. This updates every
# cell that holds one of those placeholder strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value2  = "{{{AccountName__c}}}（以下、甲という。）と　株式会社サンプル（以下、乙という。）は、"
$ws.Range("C7").Value2  = "{{{AccountName__c}}}"
$ws.Range("C8").Value2  = "{{{AccountAddress__c}}}"
$ws.Range("C11").Value2 = "{{{StartDateFormat__c}}} 〜 {{{EndDateFormat__c}}}  "
$ws.Range("C12").Value2 = "{{{Address__c}}}"
$ws.Range("C13").Value2 = "{{{JobDescription__c}}} "
$ws.Range("C14").Value2 = "{{{StartTime__c}}} 〜 {{{EndTime__c}}}  "
$ws.Range("C15").Value2 = "{{{hasOverTime__c}}} "
$ws.Range("C16").Value2 = "{{{HoliDayType__c}}} "
$ws.Range("C17").Value2 = "基本給(月)　{{{Salary__c}}}万円"
$ws.Range("C18").Value2 = "{{{DueDate__c}}} "
$ws.Range("C19").Value2 = "{{{SalaryDate__c}}} "

# Restore the selected cell to B4 (matches the author's saved view state).
$ws.Range("B4").Select() | Out-Null

# Best-effort: also try to match the saved window size recorded in the
# author's workbook view (cosmetic only; not all hosts expose this).
try { $excel.ActiveWindow.Height = 14640 } catch {}
